# Murugan Ashwin.xlsx — add matchNo column + new match row, rename sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet to the batter's name.
$ws.Name = "Murugan Ashwin"

# 2) Insert a new first column ("matchNo") shifting the existing
#    teamName..result columns from A:L to B:M.
$ws.Columns.Item(1).Insert()

# Helper: write a value as TEXT (matches source data where numeric-looking
# strings like "6", "14", "42.85" are stored as text, per the
# numberStoredAsText ignoredError already on the sheet). Temporarily force
# a text number-format so Excel doesn't coerce the digits into a real
# number, then drop the format again so no stray cell style is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# 3) Header row.
Set-TextValue $ws.Range("A1") "matchNo"

# 4) Existing match (row 2) gets its matchNo.
Set-TextValue $ws.Range("A2") "8th"

# Re-assert text storage for the numeric-looking values already in row 2
# (values themselves are unchanged, only the new column/row are additive).
Set-TextValue $ws.Range("E2") "6"
Set-TextValue $ws.Range("F2") "14"
Set-TextValue $ws.Range("G2") "0"
Set-TextValue $ws.Range("H2") "0"
Set-TextValue $ws.Range("I2") "42.85"

# 5) New match row (row 3).
Set-TextValue $ws.Range("A3") "14th"
$ws.Range("B3").Value = "Punjab Kings"
$ws.Range("C3").Value = "Murugan Ashwin"
$ws.Range("D3").Value = "c " + [char]0x2020 + "Bairstow b Kaul"
Set-TextValue $ws.Range("E3") "9"
Set-TextValue $ws.Range("F3") "10"
Set-TextValue $ws.Range("G3") "1"
Set-TextValue $ws.Range("H3") "0"
Set-TextValue $ws.Range("I3") "90.00"
$ws.Range("J3").Value = "Sunrisers Hyderabad"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 21"
$ws.Range("M3").Value = "Sunrisers won by 9 wickets (with 8 balls remaining)"
